$wb = $excel.ActiveWorkbook

# --- 1. Insert a new header-style row into "Univers" right after the title row ---
$univers = $wb.Worksheets.Item("Univers")
$univers.Rows("2:2").Insert()
$univers.Range("A2").Value = "Ma question"
$univers.Range("B2").Value = "Ma réponse"

# Return the view to the top of the sheet (selection A1:B1, no scrolled topLeftCell)
$univers.Activate()
$univers.Range("A1:B1").Select()

# --- 2. Insert a brand-new worksheet "Nouveau truc" right after "Univers" ---
$newSheet = $wb.Worksheets.Add($null, $univers)
$newSheet.Name = "Nouveau truc"

# Match the look of the other FAQ sheets: copy the header formatting from Univers!A1:B1
$univers.Range("A1:B1").Copy()
$newSheet.Range("A1:B1").PasteSpecial(-4122)
$newSheet.Rows("1:1").RowHeight = 30

$newSheet.Range("A1").Value = "Sujets"
$newSheet.Range("B1").Value = "Observations"
$newSheet.Range("A2").Value = "Biduel"
$newSheet.Range("B2").Value = "machin"

# Make the new sheet the active tab, with A3 selected
$newSheet.Activate()
$newSheet.Range("A3").Select()

Write-Output "done"
